$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column I entirely: the updated program only spans 8 weeks (A-H)
# instead of the previous 9 (A-I), so the whole last column is dropped.
$ws.Columns("I").Delete()

# Update cell values per the revised weekly program schedule (Jan-Feb 2025).
# The header row holds dd/mm/yyyy dates as plain text; prefix with a leading
# apostrophe (just like typing it in Excel) so day/month-ambiguous values such
# as "07/01/2025" are not auto-converted into date serial numbers.
$ws.Range('A1').Value = "'07/01/2025"
$ws.Range('B1').Value = "'14/01/2025"
$ws.Range('C1').Value = "'21/01/2025"
$ws.Range('D1').Value = "'28/01/2025"
$ws.Range('E1').Value = "'04/02/2025"
$ws.Range('F1').Value = "'11/02/2025"
$ws.Range('G1').Value = "'18/02/2025"
$ws.Range('H1').Value = "'25/02/2025"
$ws.Range('A2').Value = 'SALMOS 127-134'
$ws.Range('B2').Value = 'SALMOS 135-137'
$ws.Range('C2').Value = 'SALMOS 138,139'
$ws.Range('D2').Value = 'SALMOS 140-143'
$ws.Range('E2').Value = 'SALMOS 144-146'
$ws.Range('F2').Value = 'SALMOS 147-150'
$ws.Range('G2').Value = 'PROVERBIOS 1'
$ws.Range('H2').Value = 'PROVERBIOS 2'
$ws.Range('A3').Value = 'Canción 134'
$ws.Range('B3').Value = 'Canción 2'
$ws.Range('C3').Value = 'Canción 93'
$ws.Range('D3').Value = 'Canción 44'
$ws.Range('E3').Value = 'Canción 145'
$ws.Range('F3').Value = 'Canción 12'
$ws.Range('G3').Value = 'Canción 88'
$ws.Range('H3').Value = 'Canción 35'
$ws.Range('A5').Value = '1. Padres, sigan cuidando la herencia que Jehová les dio'
$ws.Range('B5').Value = '1. “Nuestro Señor es más grande que todos los demás dioses”'
$ws.Range('C5').Value = '1. ¡Que los nervios no lo frenen!'
$ws.Range('D5').Value = '1. ¿Qué hará después de orar?'
$ws.Range('E5').Value = '1. “¡Feliz el pueblo que tiene por Dios a Jehová!”'
$ws.Range('F5').Value = '1. Tenemos muchas razones para alabar a Jah'
$ws.Range('G5').Value = '1. Joven, ¿a quién escucharás?'
$ws.Range('H5').Value = '1. Por qué estudiar con ganas'
$ws.Range('B9').Value = '4. Empiece conversaciones'
$ws.Range('H9').Value = '4. Empiece conversaciones'
$ws.Range('B10').Value = '5. Haga revisitas'
$ws.Range('C10').Value = '5. Haga discípulos'
$ws.Range('F10').Value = '5. Empiece conversaciones'
$ws.Range('H10').Value = '5. Haga revisitas'
$ws.Range('A11').Value = '6. Haga discípulos'
$ws.Range('B11').Value = '6. Explique sus creencias'
$ws.Range('C11').Value = '6. Discurso'
$ws.Range('D11').Value = '6. Explique sus creencias'
$ws.Range('E11').Value = '6. Discurso'
$ws.Range('H11').Value = '6. Discurso'
$ws.Range('A14').Value = 'Canción 13'
$ws.Range('B14').Value = 'Canción 10'
$ws.Range('C14').Value = 'Canción 59'
$ws.Range('D14').Value = 'Canción 141'
$ws.Range('E14').Value = 'Canción 59'
$ws.Range('F14').Value = 'Canción 159'
$ws.Range('G14').Value = 'Canción 89'
$ws.Range('H14').Value = 'Canción 96'
$ws.Range('A15').Value = '7. Padres, ¿están usando esta herramienta tan potente?'
$ws.Range('B15').Value = '7. Necesidades de la congregación'
$ws.Range('C15').Value = '7. Aunque sea tímido, puede tener éxito sirviendo a Jehová'
$ws.Range('D15').Value = '7. Esté preparado por si necesita atención médica o una intervención quirúrgica'
$ws.Range('E15').Value = '7. Jehová quiere que usted sea feliz'
$ws.Range('F15').Value = '7. Informe de servicio anual'
$ws.Range('G15').Value = '8. Necesidades de la congregación'
$ws.Range('H15').Value = '7. ¿Eres un cazatesoros?'
$ws.Range('A16').Value = 'Gracias a su ejemplo, sus hijos pueden aprender a...'
$ws.Range('B16').Value = '8. Estudio bíblico de la congregación'
$ws.Range('D16').Value = '¿ESTÁ PREPARADO?'
$ws.Range('E16').Value = '8. Necesidades de la congregación'
$ws.Range('G16').Value = '9. Estudio bíblico de la congregación'
$ws.Range('H16').Value = 'SUGERENCIA'
$ws.Range('A17').Value = '8. Estudio bíblico de la congregación'
$ws.Range('B17').Value = 'Palabras de conclusión(3 mins.)|Canción 90y oración'
$ws.Range('C17').Value = 'Palabras de conclusión(3 mins.)|Canción 151y oración'
$ws.Range('D17').Value = '8. Estudio bíblico de la congregación'
$ws.Range('E17').Value = '9. Estudio bíblico de la congregación'
$ws.Range('F17').Value = 'Palabras de conclusión(3 mins.)|Canción 37y oración'
$ws.Range('G17').Value = 'Palabras de conclusión(3 mins.)|Canción 80y oración'
$ws.Range('H17').Value = '8. Estudio bíblico de la congregación'
$ws.Range('A18').Value = 'Palabras de conclusión(3 mins.)|Canción 73y oración'
$ws.Range('D18').Value = 'Palabras de conclusión(3 mins.)|Canción 103y oración'
$ws.Range('E18').Value = 'Palabras de conclusión(3 mins.)|Canción 85y oración'
$ws.Range('H18').Value = 'Palabras de conclusión(3 mins.)|Canción 102y oración'

# Clear cells that no longer have content in the updated program.
$ws.Range('A12').ClearContents()
